$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: keep the F1 World Champion question (previously row 5), dropping the Innsbruck/Location one
$ws.Range("A2").Value = "Who was the F1 World Champion in 2022?"
$ws.Range("B2").Value = "Max Verstappen"
$ws.Range("C2").Value = "Person"

# Row 3: replace the Barcelona question with a new one about the song Thriller
$ws.Range("A3").Value = "Who is the artist behind the song Thriller?"
$ws.Range("B3").Value = "Michael Jackson"
$ws.Range("C3").Value = "Person"

# Row 4: replace the Italian-language question with one about the current US president
$ws.Range("A4").Value = "Who is the current president of the US?"
$ws.Range("B4").Value = "Joe Biden"
$ws.Range("C4").Value = "Person"

# The remaining old rows (5-10: mayor of Innsbruck, Facebook founder, Schumacher titles,
# Chelsea Champions League) are no longer needed - remove them entirely.
$ws.Rows("5:10").Delete()

# Restore the selection to where the user last clicked before saving.
$ws.Range("L9").Select() | Out-Null
